$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AS2").Value = "2602"
$ws.Range("AR3").Value = 20
$ws.Range("AR3").Interior.Color = 16777215
$ws.Range("AS3").Value = "4489"
$ws.Range("AR4").Value = 35
$ws.Range("AR4").Interior.Color = 32768
$ws.Range("AS4").Value = "4641"
$ws.Range("AS5").Value = "4353"
$ws.Range("AS6").Value = "4400"
$ws.Range("AS7").Value = "4689"
$ws.Range("AS8").Value = "4128"
$ws.Range("AS9").Value = "4260"
$ws.Range("AS10").Value = "4278"
$ws.Range("AS11").Value = "3865"
$ws.Range("AS13").Value = "4881"
$ws.Range("AS14").Value = "4304"
$ws.Range("AS16").Value = "4415"
$ws.Range("AR17").Value = 33
$ws.Range("AR17").Interior.Color = 32768
$ws.Range("AS17").Value = "3788"
$ws.Range("AS18").Value = "3801"
$ws.Range("AR19").Value = 31
$ws.Range("AR19").Interior.Color = 32768
$ws.Range("AS19").Value = "4425"
$ws.Range("AS20").Value = "4916"
$ws.Range("AS21").Value = "3108"
$ws.Range("AS22").Value = "4298"
$ws.Range("AS24").Value = "4362"
$ws.Range("AS26").Value = "4488"
$ws.Range("AS27").Value = "4307"
$ws.Range("AS28").Value = "3603"
$ws.Range("E28").Value = "一馆"
$ws.Range("AS29").Value = "5051"
$ws.Range("AS30").Value = "3822"
$ws.Range("AS31").Value = "4784"
$ws.Range("AS32").Value = "4654"
$ws.Range("AR33").Value = 14
$ws.Range("AR33").Interior.Color = 65535
$ws.Range("AS33").Value = "3971"
$ws.Range("AS34").Value = "4247"
$ws.Range("AS35").Value = "4128"
$ws.Range("AS36").Value = "4200"
$ws.Range("AS37").Value = "3732"
$ws.Range("AS38").Value = "4243"
$ws.Range("AS39").Value = "4277"
$ws.Range("AS40").Value = "4335"
$ws.Range("AS42").Value = "4725"
$ws.Range("AS43").Value = "4508"
$ws.Range("AS44").Value = "4239"
$ws.Range("AS45").Value = "3998"
$ws.Range("AS46").Value = "3875"
$ws.Range("AS47").Value = "3918"
$ws.Range("AS48").Value = "4573"
$ws.Range("AS49").Value = "3273"
$ws.Range("AS50").Value = "4617"
$ws.Range("AS51").Value = "2548"
$ws.Range("AS52").Value = "2552"
$ws.Range("AR54").Value = 34
$ws.Range("AR54").Interior.Color = 32768
$ws.Range("AS54").Value = "4488"
$ws.Range("AS55").Value = "3576"
$ws.Range("AS56").Value = "3625"
$ws.Range("AS57").Value = "4166"
$ws.Range("AS58").Value = "3689"
$ws.Range("AS61").Value = "4072"
$ws.Range("AS63").Value = "4153"
$ws.Range("AR64").Value = 7
$ws.Range("AR64").Interior.Color = 65535
$ws.Range("AS64").Value = "2669"
$ws.Range("AR66").Value = 6
$ws.Range("AR66").Interior.Color = 65535
$ws.Range("AS66").Value = "2580"
$ws.Range("AS67").Value = "4027"
$ws.Range("AS68").Value = "3990"
$ws.Range("AS69").Value = "3645"
$ws.Range("AR70").Value = 1
$ws.Range("AR70").Interior.Color = 65535
$ws.Range("AS70").Value = "2512"
$ws.Range("AR71").Value = 23
$ws.Range("AR71").Interior.Color = 16777215
$ws.Range("AS71").Value = "3781"
$ws.Range("AS72").Value = "3159"
$ws.Range("AS73").Value = "3655"
$ws.Range("AS74").Value = "2705"
$ws.Range("AS75").Value = "3493"
$ws.Range("AS76").Value = "4191"
$ws.Range("AS78").Value = "4090"
$ws.Range("AS80").Value = "3893"
$ws.Range("AS81").Value = "4088"
$ws.Range("AS82").Value = "3871"
$ws.Range("AS83").Value = "3637"
$ws.Range("AS85").Value = "3371"
$ws.Range("AS86").Value = "2560"
$ws.Range("AS87").Value = "3823"
$ws.Range("AR88").Value = 10
$ws.Range("AR88").Interior.Color = 65535
$ws.Range("AS88").Value = "2617"
$ws.Range("AS89").Value = "2963"
$ws.Range("AR90").Value = 0
$ws.Range("AR90").Interior.Color = 255
$ws.Range("AS90").Value = "2602"
$ws.Range("AS91").Value = "2732"
$ws.Range("AS93").Value = "2295"
$ws.Range("AR95").Value = 3
$ws.Range("AR95").Interior.Color = 65535
$ws.Range("AS95").Value = "2422"
$ws.Range("AR97").Value = 20
$ws.Range("AR97").Interior.Color = 16777215
$ws.Range("AS97").Value = "2809"
$ws.Range("AS99").Value = "3214"
$ws.Range("AS100").Value = "2416"
$ws.Range("AS101").Value = "3793"
$ws.Range("AS102").Value = "3405"
$ws.Range("AR103").Value = 20
$ws.Range("AR103").Interior.Color = 16777215
$ws.Range("AS103").Value = "2787"
$ws.Range("AR104").Value = 31
$ws.Range("AR104").Interior.Color = 32768
$ws.Range("AS104").Value = "3864"
$ws.Range("AR105").Value = 17
$ws.Range("AR105").Interior.Color = 65535
$ws.Range("AS105").Value = "3495"
$ws.Range("AS107").Value = "2557"
$ws.Range("AR109").Value = 18
$ws.Range("AR109").Interior.Color = 65535
$ws.Range("AS109").Value = "3321"
$ws.Range("AS110").Value = "3265"
$ws.Range("AS111").Value = "2545"
$ws.Range("AS112").Value = "2563"
$ws.Range("AS113").Value = "3093"
$ws.Range("AS114").Value = "3155"
$ws.Range("AS115").Value = "2020"
$ws.Range("AS116").Value = "2917"
$ws.Range("AS117").Value = "3259"
$ws.Range("AR119").Value = 16
$ws.Range("AR119").Interior.Color = 65535
$ws.Range("AS119").Value = "3483"
$ws.Range("AR120").Value = 5
$ws.Range("AR120").Interior.Color = 65535
$ws.Range("AS120").Value = "2481"
$ws.Range("AR121").Value = 0
$ws.Range("AR121").Interior.Color = 255
$ws.Range("AS121").Value = "2947"
$ws.Range("AS123").Value = "2845"
$ws.Range("AS124").Value = "2995"
$ws.Range("AR125").Value = 20
$ws.Range("AR125").Interior.Color = 16777215
$ws.Range("AS125").Value = "2904"
$ws.Range("AS126").Value = "3048"
$ws.Range("AS127").Value = "2339"
$ws.Range("AS130").Value = "2789"
$ws.Range("AS132").Value = "2528"
$ws.Range("AR133").Value = 0
$ws.Range("AR133").Interior.Color = 255
$ws.Range("AS133").Value = "2597"
$ws.Range("AS139").Value = "2036"
$ws.Range("AS140").Value = "2716"
$ws.Range("AS149").Value = "3342"
$ws.Range("AS151").Value = "2550"
$ws.Range("AS155").Value = "2614"
$ws.Range("AR156").Value = 10
$ws.Range("AR156").Interior.Color = 65535
$ws.Range("AS156").Value = "3190"
$ws.Range("AS159").Value = "2821"
$ws.Range("AS161").Value = "2988"
$ws.Range("AR162").Value = 20
$ws.Range("AR162").Interior.Color = 16777215
$ws.Range("AS162").Value = "2656"
$ws.Range("AS163").Value = "3219"
$ws.Range("AS164").Value = "3016"
$ws.Range("AS165").Value = "2759"
$ws.Range("AS166").Value = "2768"
$ws.Range("AR170").Value = 3
$ws.Range("AR170").Interior.Color = 65535
$ws.Range("AS170").Value = "1517"
$ws.Range("AR171").Value = 0
$ws.Range("AR171").Interior.Color = 255
$ws.Range("AS173").Value = "1758"
